$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Devops tools")

# Fix typo "Coninuum" -> "Continuum" in the CI tools list (cell B3)
$ws.Range("B3").Value = "ElectricFlow/Cloud, Jenkins, Hudson,Bamboo, CodeShip, Travis CI, TeamCity,Continuum, BuildMaster, QuickBuild,Snap CI,CircleCI, CruiseControl,Gump,Shippable,Urban Code Build,Continua CI,Visual Studio/TFS,IBM Tivoli,Solano CI,Websphere Commerce Server,...."

# Move the active selection to reflect where the user ended up after editing
$ws.Activate()
$ws.Range("B4").Select()
